# Updating project plan with current effort's status.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Mark current task statuses in the Status column (B) ---
# Most newly-started/older tasks are now "Done"; a few active ones are "WIP".
$doneRows = @(8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 21, 22)
foreach ($r in $doneRows) {
    $ws1.Range("B$r").Value = "Done"
}

$wipRows = @(20, 23, 24)
foreach ($r in $wipRows) {
    $ws1.Range("B$r").Value = "WIP"
}

# --- Drop the stale hidden Gantt chart-tracking defined names ---
while ($wb.Names.Count() -gt 0) {
    $wb.Names.Item(1).Delete()
}

# --- Update sheet selection / active tab: Sheet1 is now the active sheet ---
$ws1.Activate()
$ws1.Range("D25").Select()
